$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Section 1 (rows 32-41): Left/Right example
#   E32 = individual formula, E33:E41 share one formula group
# ---------------------------------------------------------------
$ws.Range("E32").Formula = '=IF(H32="Left",LEFT(G32,1),RIGHT(G32,1))'
$ws.Range("E33:E41").Formula = '=IF(H33="Left",LEFT(G33,1),RIGHT(G33,1))'

# ---------------------------------------------------------------
# Section 2 (rows 57-76): MID example
#   E57 = individual formula, E58:E76 share one formula group
# ---------------------------------------------------------------
$ws.Range("E57").Formula = '=MID(G57,H57,I57)'
$ws.Range("E58:E76").Formula = '=MID(G58,H58,I58)'

# ---------------------------------------------------------------
# New reference link cell, row 89 (column G), style copied from C89
# ---------------------------------------------------------------
$ws.Range("C89").Copy()
$ws.Range("G89").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G89").Value = "https://www.youtube.com/watch?v=d-NxycY_fYQ"

# ---------------------------------------------------------------
# Section 3 (rows 98-117): INDEX/TEXTSPLIT example, each row is its
# own (non-shared) array-entered formula
# ---------------------------------------------------------------
for ($r = 98; $r -le 117; $r++) {
    $ws.Range("E$r").FormulaArray = "=INDEX(TEXTSPLIT(G$r,`" `"),H$r)"
}

# ---------------------------------------------------------------
# Section 4 (rows 138-157): SWITCH/LOWER/PROPER/UPPER example, each
# row is its own (non-shared) array-entered formula
# ---------------------------------------------------------------
for ($r = 138; $r -le 157; $r++) {
    $ws.Range("E$r").FormulaArray = "=SWITCH(H$r, `"Lower`",LOWER(G$r),`"Proper`",PROPER(G$r),`"Upper`",UPPER(G$r))"
}

# ---------------------------------------------------------------
# Update the view: scroll down and move the active selection to E138
# ---------------------------------------------------------------
$ws.Range("E138").Select() | Out-Null
